$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Detain, Remove Individual Protesters, Force: Vague/Body, Monitor/Present"
$ws.Range("D3").Value = "Monitor/Present, Remove Individual Protesters"
$ws.Range("D4").Value = "Monitor/Present, Constrain, Arrest or Attempted"
$ws.Range("D5").Value = "Constrain, Monitor/Present"
$ws.Range("D6").Value = "Constrain, Monitor/Present"
$ws.Range("D7").Value = "Force: Vague/Body, Monitor/Present"
$ws.Range("D8").Value = "Monitor/Present, Force: Vague/Body"
$ws.Range("D9").Value = "Monitor/Present, Constrain, Force: Vague/Body"
$ws.Range("D10").Value = "Monitor/Present, Detain"
$ws.Range("D11").Value = "Monitor/Present, Arrest or Attempted"
$ws.Range("D12").Value = "Constrain, Instruct/Warn, Monitor/Present"
$ws.Range("D13").Value = "Instruct/Warn, Force: 2+ Weapon Types, End Protest, `"Breaking the Rules`", Constrain, Arrest or Attempted, Monitor/Present"
$ws.Range("D14").Value = "Arrest or Attempted, Constrain, Monitor/Present"
$ws.Range("D15").Value = "Instruct/Warn, Monitor/Present"
$ws.Range("D16").Value = "Constrain, Monitor/Present"
$ws.Range("D17").Value = "Monitor/Present, Instruct/Warn"
$ws.Range("D18").Value = "Arrest or Attempted, Constrain, Monitor/Present, Instruct/Warn"
$ws.Range("D19").Value = "Remove Individual Protesters, Arrest- Large Scale, Instruct/Warn, Constrain, Force: 2+ Weapon Types, Detain, Arrest or Attempted, Monitor/Present, Formal Accusation"
$ws.Range("D20").Value = "Monitor/Present, Instruct/Warn, Force: Vague/Body, Force: Weapon"
$ws.Range("D21").Value = "Arrest or Attempted, Formal Accusation, Instruct/Warn, Monitor/Present"
$ws.Range("D22").Value = "Monitor/Present, Constrain"
$ws.Range("D23").Value = "Monitor/Present, Constrain"
$ws.Range("D24").Value = "Monitor/Present, Constrain"
$ws.Range("D25").Value = "Constrain, Detain, Cooperate/Coordinate, Monitor/Present"
$ws.Range("D26").Value = "Monitor/Present, Constrain, Detain, Formal Accusation, Arrest or Attempted"
$ws.Range("D27").Value = "Monitor/Present, Remove Individual Protesters, Constrain, Detain, Formal Accusation"
$ws.Range("D28").Value = "Arrest or Attempted, Formal Accusation, Instruct/Warn, Constrain, Monitor/Present"
$ws.Range("D29").Value = "Arrest or Attempted, Formal Accusation, Constrain, Monitor/Present"
$ws.Range("D30").Value = "Monitor/Present, Constrain"
$ws.Range("D31").Value = "Detain, Monitor/Present, Arrest or Attempted"
$ws.Range("D32").Value = "Monitor/Present, Instruct/Warn"
$ws.Range("D33").Value = "Formal Accusation, `"Breaking the Rules`", Monitor/Present, Arrest or Attempted"
$ws.Range("D34").Value = "Monitor/Present, Constrain"
$ws.Range("D35").Value = "Monitor/Present, Instruct/Warn"
$ws.Range("D36").Value = "Monitor/Present, Instruct/Warn, `"Breaking the Rules`""
$ws.Range("D37").Value = "Instruct/Warn, Monitor/Present, `"Breaking the Rules`", Detain, Formal Accusation, End Protest"
$ws.Range("D38").Value = "Monitor/Present, Formal Accusation, Arrest or Attempted, End Protest"
$ws.Range("D39").Value = "Arrest or Attempted, Formal Accusation, Monitor/Present"
$ws.Range("D40").Value = "Monitor/Present, Instruct/Warn"
